$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the "Innsbruck" question (row 2) and before the
# Schumacher question, to hold the new Tyrol question/answer pair.
$ws.Rows.Item(3).Insert()

# Row 3: new Tyrol question
$ws.Range("A3").Value = "What is the capital of Tyrol?"
$ws.Range("B3").Value = "Innsbruck"
$ws.Range("C3").Value = "Location"

# Row 5: new "last title" Schumacher question
$ws.Range("A5").Value = "When did Michael Schumacher win his last F1 World Drivers Title?"
$ws.Range("B5").Value = 2004
$ws.Range("C5").Value = "Year"

# Row 6: Grand Tour host question (Richard Hammond)
$ws.Range("A6").Value = "Who hosts The Grand Tour?"
$ws.Range("B6").Value = "Richard Hammond"
$ws.Range("C6").Value = "Person"

# Row 7: reword the 2022 F1 champion question, keep Max Verstappen answer
$ws.Range("A7").Value = "Who was 2022 F1 World Drivers champion?"
$ws.Range("B7").Value = "Max Verstappen"
$ws.Range("C7").Value = "Person"

# Row 8: Grand Tour host question (Jeremy Clarkson)
$ws.Range("A8").Value = "Who hosts The Grand Tour?"
$ws.Range("B8").Value = "Jeremy Clarkson"
$ws.Range("C8").Value = "Person"

# Widen column B slightly to fit the new longer values
$ws.Columns.Item(2).ColumnWidth = 16.42

# Page setup (paper size / orientation) as set when the file was last printed
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the selected cell as recorded in the saved workbook
$ws.Range("L13").Select()
